$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2: swap order of the product summary ("2-queque,1-torta," -> "1-torta,2-queque,")
$ws.Range("A2").Value = "1-torta,2-queque,"

# Delete rows 7-10 (test/dummy orders), shifting rows up
$ws.Range("A7:J10").EntireRow.Delete()
